$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 350000
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 450000
$ws.Range("J2").Value = 45775
$ws.Range("K2").Value = 45805
$ws.Range("M2").Value = 5
